$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the formatting (border style) from the row above onto the two new rows
$ws.Range("A18:C18").Copy()
$ws.Range("A19:C20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new values (column-by-column, matching the order the
# author originally typed them in)
$ws.Range("A19").Value = "sex"
$ws.Range("A20").Value = "age"
$ws.Range("B19").Value = "Gender"
$ws.Range("B20").Value = "Age"
$ws.Range("C19").Value = "World Bank"
$ws.Range("C20").Value = "World Bank"

# Update the current selection to match the author's final state
$ws.Range("A3:C20").Select()
